# Add a new worksheet "Event Type Name List" at the end of the workbook,
# containing the accepted list of event type names (one per row in column A),
# with a bold/centered/bordered header in A1 matching the existing
# "Appointment Type" header style used on the "Appointment Type Summation" sheet.

$wb = $excel.ActiveWorkbook

# Sheet used as the source of the header formatting (existing header style).
$styleSourceSheet = $wb.Worksheets.Item("Appointment Type Summation")

# Add the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Event Type Name List"

$values = @(
    "Event Type Name Accepted List",
    "Appointment",
    "Big Interview ",
    "Career Closet",
    "Career Course",
    "Career Fair",
    "Classroom Presentation",
    "Club Support ",
    "Club Presentation ",
    "Completed Handshake Profile",
    "Drop-In/Chat",
    "Employer Partner Event",
    "Employment Toolkit",
    "Hiration",
    "HS Employer Review",
    "HS Interview Review",
    "Info Session",
    "Library Book",
    "Mentor Meetup ",
    "Networking",
    "Other",
    "Possible Program (Fall Only?)",
    "Project Onramp",
    "Rise Together",
    "Speaker/Panel",
    "Trek",
    "Type Focus",
    "Workshop",
    "WOW"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $values[$i]
}

# Match the bold/centered/top-aligned/bordered header style already used
# for header cells on the other sheets (e.g. "Appointment Type" on the
# "Appointment Type Summation" sheet).
$styleSourceSheet.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
